$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Products sheet: add a "Quantity" column (C) next to the existing
# Products/Cost table.
# ---------------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("Products")

# Give C1 the same look (border / number format / alignment) as the existing
# "Cost" header in B1 before we put the new header text in it.
$wsProducts.Range("B1").Copy()
$wsProducts.Range("C1").PasteSpecial(-4122)
$wsProducts.Range("C1").Value = "Quantity"

$wsProducts.Range("C2").Value = 120
$wsProducts.Range("C3").Value = 100
$wsProducts.Range("C4").Value = 50

# Bold the whole header row (Products / Cost / Quantity).
$wsProducts.Range("A1:C1").Font.Bold = $true

# ---------------------------------------------------------------------------
# Expenses sheet: drop in a small reference table of names.
# ---------------------------------------------------------------------------
$wsExpenses = $wb.Worksheets.Item("Expenses")

$wsExpenses.Range("A1").Value = 1
$wsExpenses.Range("B1").Value = "Jason"
$wsExpenses.Range("C1").Value = "Brody"

$wsExpenses.Range("A2").Value = 2
$wsExpenses.Range("B2").Value = "John"
$wsExpenses.Range("C2").Value = "Smith"

$wsExpenses.Range("A3").Value = 3
$wsExpenses.Range("B3").Value = "Mark"
$wsExpenses.Range("C3").Value = "Hopper"

$wsExpenses.Range("A4").Value = 4
$wsExpenses.Range("B4").Value = "Janis"
$wsExpenses.Range("C4").Value = "Joplin"

$wsExpenses.Range("A5").Value = 5
$wsExpenses.Range("B5").Value = "Jimi"
$wsExpenses.Range("C5").Value = "Hendrix"

# ---------------------------------------------------------------------------
# Clients sheet: new client write functionality - ClientID / ClientName /
# Column1 table.
# ---------------------------------------------------------------------------
$wsClients = $wb.Worksheets.Item("Clients")

$wsClients.Range("B1").Value = "ClientName"
$wsClients.Range("C1").Value = "Column1"

$wsClients.Range("A2").Value = 1
$wsClients.Range("B2").Value = "James"

$wsClients.Range("A1").Value = "ClientID"

$wsClients.Range("A3").Value = 1
$wsClients.Range("B3").Value = "Semaj"

$wsClients.Range("A1:C1").Font.Bold = $true

$wsClients.Columns("A").ColumnWidth = 12.14
$wsClients.Columns("B").ColumnWidth = 11.14

# ---------------------------------------------------------------------------
# Final view state: selection was left on D4 of Products, then the user
# moved to the Clients sheet and left the selection on A2.
# ---------------------------------------------------------------------------
$wsProducts.Range("D4").Select()
$wsClients.Activate()
$wsClients.Range("A2").Select()
